$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.442.13"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "'1.844.96"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'262.84"
$ws.Range("E5").Value = "  -3.83%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D7").Value = "'0.5207"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("D8").Value = "'0.3265"
$ws.Range("E8").Value = "  -3.49%  "
$ws.Range("D9").Value = "'0.06797"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -5.50%  "
$ws.Range("D11").Value = "'0.7793"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "'0.07747"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "'1.844.82"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'88.11"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "'5.013"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("D17").Value = "'13.93"
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").Value = "'0.000007974"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "'0.9995"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'26.483.16"
$ws.Range("D21").Value = "'2.076.23"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "'4.618"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").Value = "'9.580"
$ws.Range("E23").Value = "  -3.75%  "
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("D25").Value = "'144.37"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("E26").Value = "  -7.86%  "
$ws.Range("D27").Value = "'1.654"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("D29").Value = "'112.10"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'4.159"
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("D31").Value = "'4.127"
$ws.Range("E31").Value = "  -4.15%  "
$ws.Range("D32").Value = "'0.08708"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "'0.04832"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").Value = "'0.7213"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").Value = "'1.131"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").Value = "'2.839"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").Value = "'3.111"
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01781"
$ws.Range("E38").Value = "  -3.42%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.228"
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("D40").Value = "'0.4862"
$ws.Range("E40").Value = "  -4.41%  "
$ws.Range("D41").Value = "'0.9129"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").Value = "'111.08"
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("D43").Value = "'6.070"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D44").Value = "'0.9995"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'7.739"
$ws.Range("D46").Value = "'0.4180"
$ws.Range("E46").Value = "  -5.08%  "
$ws.Range("D47").Value = "'0.05932"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "'9.064"
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1236"
$ws.Range("E49").Value = "  -6.72%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'35.06"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "'0.8883"
$ws.Range("E51").Value = "  +1.02%  "
